$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New projected values for rows 2-6, columns B through Q
# (age calculation + fixed tenure for new hires in projection utils)
$data = @{
    2 = @{ B=102; C=102; D=86;  E=0.8431372549019608;  F=0.8431372549019608;  G=0.09853558668981278; H=0.08307902407180294; I=453398.1905958019;   J=164085.0954989009; K=0; L=164085.0954989009; M=617483.2860947028;  N=10034971.8888;      O=9627230.958700001;  P=0.01635132587486724; Q=0.01704385157090465 }
    3 = @{ B=103; C=103; D=85;  E=0.8252427184466019;  F=0.8252427184466019;  G=0.09946524789627358; H=0.08208297156488596; I=475130.6665414795;   J=172352.9134777698; K=0; L=172352.9134777698; M=647483.5800192493;  N=10633646.086764;    O=10226272.928761;    P=0.01620826121835127; Q=0.01685393248140618 }
    4 = @{ B=104; C=104; D=87;  E=0.8365384615384616;  F=0.8365384615384616;  G=0.09774377642832541; H=0.08176642835831067; I=502966.1150419703;   J=179184.265134226;  K=0; L=179184.265134226;  M=682150.3801761963;  N=10970666.81506692;  O=10562322.46232383;  P=0.01633303318337382; Q=0.01696447592595118 }
    5 = @{ B=105; C=104; D=87;  E=0.8365384615384616;  F=0.8285714285714286;  G=0.09774030428156209; H=0.08098482354758002; I=516886.2289629109;   J=183570.3583230936; K=0; L=183570.3583230936; M=700456.5872860046;  N=11269370.82291893;  O=10858676.13959355;  P=0.01628931740801003; Q=0.01690540872231639 }
    6 = @{ B=106; C=106; D=90;  E=0.8490566037735849;  F=0.8490566037735849;  G=0.09459020222215261; H=0.08031243584899751; I=533100.1538977289;   J=188616.1381057517; K=0; L=188616.1381057517; M=721716.2920034805;  N=11647629.9738065;   O=11233164.44998135;  P=0.01619352078748352; Q=0.0167910065721565 }
}

$columns = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    foreach ($col in $columns) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}
